$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Update Row 6: Orchestrator Credential entry (previously empty row)
# Cell entry order matches shared-string insertion order captured in the diff.
$ws.Range("B6").Value = "OrchestratorCredential"
$ws.Range("A6").Value = "Orchestrator Credential"
$ws.Range("C6").Value = "Orchestrator platform username and password."

# Update Row 3: One Drive folder name / value (description cell reuses existing string)
$ws.Range("A3").Value = "One Drive Folder Name"
$ws.Range("B3").Value = "OneDrive - PRP Business Solutions"
$ws.Range("C3").Value = "OneDrive folder location where Faulty Bills will be saved. "

# Update selection/active cell to A3
$ws.Range("A3").Select()

$wb.Save()
